$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing values
$ws.Range("B2").Value = 4.5
$ws.Range("B3").Value = 4.5
$ws.Range("C5").Value = 25

# Add a new value in K1
$ws.Range("K1").Value = 1

# Move selection to C3 as in the final saved state
$ws.Range("C3").Select()
